$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# Build the two new border-only styles once, on sheet 1's C1 / D1:
#   C1 -> top+bottom border only
#   D1 -> top+bottom+right border only
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# Clone the two styles built above onto C1/D1 and F1/G1 instead of
# rebuilding them with border ops again, so every cell reuses the same
# two new style records rather than minting extra ones.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 becomes genuinely empty (drop the stray empty inline-string cell)
$ws2.Range("G5").ClearContents()
